# Scheduled market-data refresh: updates price/profit columns (H-N)
# on the Leve-profit tables across all crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28: The Writing Is Not on the Wall
$ws.Range("H28").Value = 2511.6155
$ws.Range("I28").Value = 468.1
$ws.Range("J28").Value = 9323.333000000001
$ws.Range("K28").Value = 468.1
$ws.Range("L28").Value = 9323.333000000001
$ws.Range("M28").Value = 16.89999999999998
$ws.Range("N28").Value = -10293.333

# Row 33: Glazed and Confused
$ws.Range("H33").Value = 197.72223
$ws.Range("I33").Value = 213.93333
$ws.Range("J33").Value = 116.666664
$ws.Range("K33").Value = 213.93333
$ws.Range("L33").Value = 116.666664
$ws.Range("M33").Value = 15.06666999999999
$ws.Range("N33").Value = -574.666664

# Row 88: The Grave of Hemlock Groves
$ws.Range("H88").Value = 2096.5454
$ws.Range("J88").Value = 2624.875
$ws.Range("L88").Value = 2624.875
$ws.Range("N88").Value = -3436.875

# Row 91: Dappling the Highlands (L)
$ws.Range("H91").Value = 2096.5454
$ws.Range("J91").Value = 2624.875
$ws.Range("L91").Value = 2624.875
$ws.Range("N91").Value = -5432.875

# Row 107: Another Man's Ink
$ws.Range("H107").Value = 372.81818
$ws.Range("I107").Value = 282.77777
$ws.Range("K107").Value = 282.77777
$ws.Range("M107").Value = 1637.22223

# Row 129: Practical Command
$ws.Range("H129").Value = 1092.1613
$ws.Range("J129").Value = 1161.8572
$ws.Range("L129").Value = 3485.5716
$ws.Range("N129").Value = -13485.5716

$ws = $wb.Worksheets.Item("ARM")
# Row 88: The Mast Chance
$ws.Range("H88").Value = 2815
$ws.Range("J88").Value = 2918.75
$ws.Range("L88").Value = 2918.75
$ws.Range("N88").Value = -3730.75

# Row 91: The Rose and the Riveter (L)
$ws.Range("H91").Value = 2815
$ws.Range("J91").Value = 2918.75
$ws.Range("L91").Value = 2918.75
$ws.Range("N91").Value = -5726.75

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 876.8889
$ws.Range("I110").Value = 879.8
$ws.Range("J110").Value = 873.25
$ws.Range("K110").Value = 879.8
$ws.Range("L110").Value = 873.25
$ws.Range("M110").Value = 1165.2
$ws.Range("N110").Value = -4963.25

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 5121.75
$ws.Range("I122").Value = 7030.143
$ws.Range("J122").Value = 2450
$ws.Range("K122").Value = 21090.429
$ws.Range("L122").Value = 7350
$ws.Range("M122").Value = -18640.429
$ws.Range("N122").Value = -12250

# Row 131: Additions to the Armoire
$ws.Range("H131").Value = 67428.75
$ws.Range("J131").Value = 67428.75
$ws.Range("L131").Value = 67428.75
$ws.Range("N131").Value = -77508.75

$ws = $wb.Worksheets.Item("BSM")
# Row 64: With Bearings Straight
$ws.Range("H64").Value = 686.8889
$ws.Range("I64").Value = 565.8333
$ws.Range("J64").Value = 929
$ws.Range("K64").Value = 565.8333
$ws.Range("L64").Value = 929
$ws.Range("M64").Value = -340.8333
$ws.Range("N64").Value = -1379

# Row 67: Bearing the Brunt (L)
$ws.Range("H67").Value = 686.8889
$ws.Range("I67").Value = 565.8333
$ws.Range("J67").Value = 929
$ws.Range("K67").Value = 565.8333
$ws.Range("L67").Value = 929
$ws.Range("M67").Value = 214.1667
$ws.Range("N67").Value = -2489

# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 128459.625
$ws.Range("I86").Value = 4701.25
$ws.Range("J86").Value = 252218
$ws.Range("K86").Value = 4701.25
$ws.Range("L86").Value = 252218
$ws.Range("M86").Value = -3578.25
$ws.Range("N86").Value = -254464

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 128459.625
$ws.Range("I89").Value = 4701.25
$ws.Range("J89").Value = 252218
$ws.Range("K89").Value = 23506.25
$ws.Range("L89").Value = 1261090
$ws.Range("M89").Value = -17890.25
$ws.Range("N89").Value = -1272322

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2432.3784
$ws.Range("I134").Value = 2275.7932
$ws.Range("K134").Value = 6827.3796
$ws.Range("M134").Value = -4292.3796

# Row 138: Bladewinner
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# Row 140: Ceremonial Teeth
$ws.Range("H140").Value = 66303.336
$ws.Range("J140").Value = 66303.336
$ws.Range("L140").Value = 66303.336
$ws.Range("N140").Value = -76663.336

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 2500
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -1713
$ws.Range("N16").Value = -3574

# Row 64: Almost as Fun as Slingshotting Birds
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 67: Living Bow to Mouth (L)
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 70: A Reward Fitting of the Faithful
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73: Just Rewards for Just Devotion (L)
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# Row 88: Hold on Adamantite
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

# Row 91: Spears for Stone Vigilantes (L)
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

# Row 113: Patient Patients
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("CUL")
# Row 122: Salt of the North
$ws.Range("H122").Value = 732
$ws.Range("I122").Value = 450.33334
$ws.Range("K122").Value = 4053.00006
$ws.Range("M122").Value = -1603.00006

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 2680.3333
$ws.Range("I102").Value = 2291.6667
$ws.Range("J102").Value = 3263.3333
$ws.Range("K102").Value = 2291.6667
$ws.Range("L102").Value = 3263.3333
$ws.Range("M102").Value = -669.6667000000002
$ws.Range("N102").Value = -6507.3333

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 4176
$ws.Range("I122").Value = 3766.5
$ws.Range("J122").Value = 4831.2
$ws.Range("K122").Value = 11299.5
$ws.Range("L122").Value = 14493.6
$ws.Range("M122").Value = -8849.5
$ws.Range("N122").Value = -19393.6

# Row 131: Star Athletes
$ws.Range("H131").Value = 37653
$ws.Range("J131").Value = 37653
$ws.Range("L131").Value = 37653
$ws.Range("N131").Value = -47733

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad
$ws.Range("H40").Value = 4660.294
$ws.Range("I40").Value = 4877.625
$ws.Range("J40").Value = 4467.1113
$ws.Range("K40").Value = 4877.625
$ws.Range("L40").Value = 4467.1113
$ws.Range("M40").Value = -4741.625
$ws.Range("N40").Value = -4739.1113

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 2803.45
$ws.Range("I61").Value = 2471.7334
$ws.Range("K61").Value = 2471.7334
$ws.Range("M61").Value = -2269.7334

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 2067.5386
$ws.Range("I68").Value = 1887.8
$ws.Range("J68").Value = 2666.6667
$ws.Range("K68").Value = 1887.8
$ws.Range("L68").Value = 2666.6667
$ws.Range("M68").Value = -1138.8
$ws.Range("N68").Value = -4164.6667

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 2067.5386
$ws.Range("I71").Value = 1887.8
$ws.Range("J71").Value = 2666.6667
$ws.Range("K71").Value = 9439
$ws.Range("L71").Value = 13333.3335
$ws.Range("M71").Value = -5695
$ws.Range("N71").Value = -20821.3335

# Row 113: Peace in Rest
$ws.Range("H113").Value = 2803.45
$ws.Range("I113").Value = 2471.7334
$ws.Range("K113").Value = 2471.7334
$ws.Range("M113").Value = -301.7334000000001

# Row 122: Hell on Leather
$ws.Range("H122").Value = 10002709
$ws.Range("I122").Value = 2679.1
$ws.Range("J122").Value = 20002738
$ws.Range("K122").Value = 8037.299999999999
$ws.Range("L122").Value = 60008214
$ws.Range("M122").Value = -5587.299999999999
$ws.Range("N122").Value = -60013114

$ws = $wb.Worksheets.Item("WVR")
# Row 104: Brimming with Confidence
$ws.Range("H104").Value = 270000
$ws.Range("J104").Value = 270000
$ws.Range("L104").Value = 270000
$ws.Range("N104").Value = -276988

# Row 123: Helping Handwear
$ws.Range("H123").Value = 21543.23
$ws.Range("J123").Value = 21543.23
$ws.Range("L123").Value = 21543.23
$ws.Range("N123").Value = -31343.23

# Row 125: Color Coated
$ws.Range("H125").Value = 60643.332
$ws.Range("J125").Value = 60643.332
$ws.Range("L125").Value = 60643.332
$ws.Range("N125").Value = -70483.33199999999

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 10701.353
$ws.Range("I126").Value = 15273.6
$ws.Range("J126").Value = 4169.5713
$ws.Range("K126").Value = 45820.8
$ws.Range("L126").Value = 12508.7139
$ws.Range("M126").Value = -43350.8
$ws.Range("N126").Value = -17448.7139

# Row 131: A Better Bottom Line
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
